$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "result on <timestamp>" column header in G1, matching the
# existing F1 pattern (a shared string is added to sharedStrings.xml and
# the cell/dimension/used-range are updated automatically on save).
$ws.Range("G1").Value = "result on16-Jul-2024-05-35-15"
